# feat: add 2022-Q4 data
#
# - "总计" gets a new summary row for 2022-Q4 (inserted above the existing
#   2022-Q2 summary row, which shifts down one row).
# - The sheet that used to be named "2022-Q2" (holding the quarterly fund
#   holdings table) is repurposed to hold the new 2022-Q4 fund holdings,
#   and a fresh sheet named "2022-Q2" is added to preserve the data that
#   used to live there.

function Set-TextValue($range, [string]$text) {
    # Force a literal text value even for numeric-looking strings (keeps
    # leading zeros like fund codes, and keeps e.g. "0.97" as text rather
    # than a coerced double), then clear the resulting quote-prefix style
    # residue so the cell is left with the default/Normal style.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计": insert the new 2022-Q4 summary row above the 2022-Q2 summary row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()

# Row 2's "A" cell needs the same bold/bordered style as the other index
# cells in column A (copy format only from A3, the shifted-down old row).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 4
$wsTotal.Cells.Item(2, 4).Value = 0.04

# The old row (now row 3) keeps its "2022-Q2" data, only its running index
# in column A advances from 0 to 1.
$wsTotal.Cells.Item(3, 1).Value = 1

# ---------------------------------------------------------------------------
# 2) Duplicate the existing "2022-Q2" sheet so its current fund-holding data
#    survives on its own tab, then reuse the original tab for the new
#    2022-Q4 fund-holding data.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($null, $wsQ2)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ2Copy = $wb.Worksheets.Item(3)

$wsQ4.Name = "2022-Q4"
$wsQ2Copy.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 3) Rebuild the data rows on the (renamed) 2022-Q4 sheet. The header row
#    (row 1) is untouched. The old sheet had 2 data rows (rows 2-3); the new
#    data needs 4 rows (rows 2-5), so drop one row and insert three fresh
#    ones, copying the format from the remaining template row.
# ---------------------------------------------------------------------------
$wsQ4.Rows.Item(3).Delete()
$wsQ4.Rows.Item(3).Insert()
$wsQ4.Rows.Item(3).Insert()
$wsQ4.Rows.Item(3).Insert()
$wsQ4.Range("A2:H2").Copy()
$wsQ4.Range("A3:H5").PasteSpecial(-4122)

# Row 2: 001219 / 上投摩根动态多因子策略混合A
$wsQ4.Cells.Item(2, 1).Value = 0
Set-TextValue $wsQ4.Cells.Item(2, 2) "001219"
Set-TextValue $wsQ4.Cells.Item(2, 3) "上投摩根动态多因子策略混合A"
Set-TextValue $wsQ4.Cells.Item(2, 4) "0.97"
Set-TextValue $wsQ4.Cells.Item(2, 5) "92.08"
Set-TextValue $wsQ4.Cells.Item(2, 6) "3.74"
Set-TextValue $wsQ4.Cells.Item(2, 7) "0.0363"
$wsQ4.Cells.Item(2, 8).Value = 9

# Row 3: 167703 / 德邦量化优选股票（LOF）C
$wsQ4.Cells.Item(3, 1).Value = 1
Set-TextValue $wsQ4.Cells.Item(3, 2) "167703"
Set-TextValue $wsQ4.Cells.Item(3, 3) "德邦量化优选股票（LOF）C"
Set-TextValue $wsQ4.Cells.Item(3, 4) "0.56"
Set-TextValue $wsQ4.Cells.Item(3, 5) "88.52"
Set-TextValue $wsQ4.Cells.Item(3, 6) "0.94"
Set-TextValue $wsQ4.Cells.Item(3, 7) "0.0053"
$wsQ4.Cells.Item(3, 8).Value = 10

# Row 4: 167702 / 德邦量化优选股票（LOF）A
$wsQ4.Cells.Item(4, 1).Value = 2
Set-TextValue $wsQ4.Cells.Item(4, 2) "167702"
Set-TextValue $wsQ4.Cells.Item(4, 3) "德邦量化优选股票（LOF）A"
Set-TextValue $wsQ4.Cells.Item(4, 4) "0.32"
Set-TextValue $wsQ4.Cells.Item(4, 5) "88.52"
Set-TextValue $wsQ4.Cells.Item(4, 6) "0.94"
Set-TextValue $wsQ4.Cells.Item(4, 7) "0.0030"
$wsQ4.Cells.Item(4, 8).Value = 10

# Row 5: 017176 / 上投摩根动态多因子策略混合C
$wsQ4.Cells.Item(5, 1).Value = 3
Set-TextValue $wsQ4.Cells.Item(5, 2) "017176"
Set-TextValue $wsQ4.Cells.Item(5, 3) "上投摩根动态多因子策略混合C"
Set-TextValue $wsQ4.Cells.Item(5, 4) "0.00"
Set-TextValue $wsQ4.Cells.Item(5, 5) "92.08"
Set-TextValue $wsQ4.Cells.Item(5, 6) "3.74"
$wsQ4.Cells.Item(5, 7).Value = 0
$wsQ4.Cells.Item(5, 8).Value = 9
